# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New K (column G) values computed for rows 2-20 (s_vals)
$kValues = @{
    2  = 1
    3  = 1
    4  = 2
    5  = 0
    6  = 0
    7  = 2
    8  = 2
    9  = 0
    10 = 1
    11 = 1
    12 = 1
    13 = 1
    14 = 0
    15 = 3
    16 = 1
    17 = 1
    18 = 2
    19 = 3
    20 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
